$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.425.69'
$ws.Range('E2').Value = '  -1.01%  '

# Row 3
$ws.Range('D3').Value = '3.770.72'
$ws.Range('E3').Value = '  -1.65%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '694.52'
$ws.Range('E5').Value = '  -1.51%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.84'
$ws.Range('E6').Value = '  -2.54%  '

# Row 7
$ws.Range('D7').Value = '3.772.33'
$ws.Range('E7').Value = '  -1.57%  '

# Row 8
$ws.Range('E8').Value = '  +0.46%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -1.39%  '

# Row 10
$ws.Range('E10').Value = '  -2.26%  '

# Row 11
$ws.Range('E11').Value = '  +1.78%  '

# Row 12
$ws.Range('E12').Value = '  +3.03%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('E13').Value = '  -3.45%  '

# Row 14
$ws.Range('E14').Value = '  -3.11%  '

# Row 15
$ws.Range('D15').Value = '4.403.37'
$ws.Range('E15').Value = '  -1.75%  '

# Row 16
$ws.Range('D16').Value = '3.751.88'
$ws.Range('E16').Value = '  -0.97%  '

# Row 17
$ws.Range('D17').Value = '70.565.35'
$ws.Range('E17').Value = '  -0.80%  '

# Row 18
$ws.Range('E18').Value = '  -0.15%  '

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.34'
$ws.Range('E19').Value = '  -0.43%  '

# Row 20
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.09'
$ws.Range('E20').Value = '  -2.04%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '512.00'
$ws.Range('E21').Value = '  +3.33%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.29'
$ws.Range('E22').Value = '  -4.03%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.708'
$ws.Range('E23').Value = '  -3.98%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.03'
$ws.Range('E24').Value = '  -2.76%  '

# Row 25
$ws.Range('E25').Value = '  -4.62%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.41'
$ws.Range('E26').Value = '  +2.51%  '

# Row 27
$ws.Range('D27').Value = '3.916.66'
$ws.Range('E27').Value = '  -1.73%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  -5.07%  '

# Row 29
$ws.Range('E29').Value = '  -0.06%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.93'
$ws.Range('E30').Value = '  -7.82%  '

# Row 31
$ws.Range('E31').Value = '  -5.59%  '

# Row 32
$ws.Range('E32').Value = '  -1.06%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.23'
$ws.Range('E33').Value = '  -2.78%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.78'
$ws.Range('E34').Value = '  -2.28%  '

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.168'
$ws.Range('E35').Value = '  -4.58%  '

# Row 36
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.16'
$ws.Range('E36').Value = '  -0.65%  '

# Row 37
$ws.Range('E37').Value = '  -0.97%  '

# Row 38
$ws.Range('D38').Value = '3.733.32'
$ws.Range('E38').Value = '  -1.67%  '

# Row 39
$ws.Range('E39').Value = '  +8.12%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0987'
$ws.Range('E40').Value = '  -3.90%  '

# Row 41
$ws.Range('E41').Value = '  -0.93%  '

# Row 42
$ws.Range('E42').Value = '  -3.96%  '

# Row 44
$ws.Range('E44').Value = '  -0.07%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.10'
$ws.Range('E45').Value = '  -6.47%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '161.94'
$ws.Range('E46').Value = '  -1.29%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '48.93'
$ws.Range('E47').Value = '  +0.03%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000298'
$ws.Range('E48').Value = '  -5.02%  '

# Row 49
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '410.07'
$ws.Range('E49').Value = '  -4.50%  '

# Row 50
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.36'
$ws.Range('E50').Value = '  -1.34%  '

# Row 51
$ws.Range('E51').Value = '  -2.53%  '
